$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values stay stored as text (matches source data
# which is numeric-looking text like "43.877.49", "0.0936", etc.)
$priceCells = @("D2","D3","D5","D6","D7","D10","D11","D12","D14","D15","D16","D17","D18", `
                "D20","D21","D23","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34", `
                "D36","D38","D39","D44","D49")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.782.23"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.294.14"
$ws.Range("E3").Value = "  -1.56%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - Solana
$ws.Range("D5").Value = "97.89"
$ws.Range("E5").Value = "  +2.78%  "

# Row 6 - BNB
$ws.Range("D6").Value = "270.49"
$ws.Range("E6").Value = "  -0.24%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  -1.16%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.62%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "45.30"
$ws.Range("E10").Value = "  +0.05%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  -1.04%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "7.90"
$ws.Range("E12").Value = "  -3.02%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.57%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "15.83"
$ws.Range("E14").Value = "  +1.51%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.636.71"
$ws.Range("E15").Value = "  -1.59%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.855"
$ws.Range("E16").Value = "  -0.48%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.294.01"
$ws.Range("E17").Value = "  -1.36%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.759.34"
$ws.Range("E18").Value = "  +0.22%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  -3.45%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "72.32"
$ws.Range("E21").Value = "  -0.19%  "

# Row 22 - ImmutableX
$ws.Range("E22").Value = "  +7.75%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "233.21"
$ws.Range("E23").Value = "  -2.74%  "

# Row 24 / Row 25 - swap InternetComputer(DFINITY) and PancakeSwap
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "2.83"
$ws.Range("E24").Value = "  +11.62%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "9.13"
$ws.Range("E25").Value = "  -2.52%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.02%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "11.25"
$ws.Range("E27").Value = "  -1.59%  "

# Row 28 - WEMIXToken
$ws.Range("D28").Value = "3.46"
$ws.Range("E28").Value = "  -1.17%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -1.71%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").Value = "38.15"
$ws.Range("E30").Value = "  -0.16%  "

# Row 31 - Monero
$ws.Range("D31").Value = "176.61"
$ws.Range("E31").Value = "  +2.38%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "21.80"
$ws.Range("E32").Value = "  -3.79%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.0894"
$ws.Range("E33").Value = "  -0.48%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "5.43"
$ws.Range("E34").Value = "  -1.02%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  +0.50%  "

# Row 36 - RenderToken
$ws.Range("D36").Value = "4.74"
$ws.Range("E36").Value = "  +8.18%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  +0.26%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.0351"
$ws.Range("E38").Value = "  -3.11%  "

# Row 39 - NEARProtocol
$ws.Range("D39").Value = "3.53"
$ws.Range("E39").Value = "  +4.28%  "

# Row 40 - Algorand
$ws.Range("E40").Value = "  +0.36%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  -1.85%  "

# Row 42 - ARBITRUM
$ws.Range("E42").Value = "  +0.45%  "

# Row 43 - Celestia
$ws.Range("E43").Value = "  +0.62%  "

# Row 44 - MultiversX
$ws.Range("D44").Value = "64.54"
$ws.Range("E44").Value = "  +3.81%  "

# Row 45 - FraxShare
$ws.Range("E45").Value = "  -3.44%  "

# Row 46 - THORChain
$ws.Range("E46").Value = "  -2.58%  "

# Row 47 - Cronos
$ws.Range("E47").Value = "  -1.16%  "

# Row 48 - TrustWalletToken
$ws.Range("E48").Value = "  +0.71%  "

# Row 49 - Aave
$ws.Range("D49").Value = "98.58"
$ws.Range("E49").Value = "  -1.88%  "

# Row 50 - WOONetwork
$ws.Range("E50").Value = "  +6.08%  "

# Row 51 - Stacks
$ws.Range("E51").Value = "  +11.46%  "
